$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-08-07 20:59:15"

# Update the timestamp column (O) for all data rows (2 through 73)
for ($r = 2; $r -le 73; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# Update productAriaLabel for row 55 (N55) to reflect "Online kein Bestand"
$ws.Range("N55").Value = "Naturaline Damen String schwarz L - Online kein Bestand 9.95 Schweizer Franken"
